$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect")

# Helper: write a value that Excel would otherwise auto-convert to a number
# (e.g. "13082", "1", "2") while keeping it stored as text/string, and
# without leaving behind any new/changed cell style. We do this by writing
# a formula that evaluates to the literal text, then collapsing the formula
# down to its computed (string) value in place via copy / paste-special
# values. This avoids the "quote prefix" / text-number-format styles that
# a leading apostrophe or a NumberFormat="@" change would otherwise add.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

# Row 2
Set-TextValue $ws.Range("A2") "13082"
$ws.Range("B2").Value = "Scenario1"
$ws.Range("C2").Value = "1) Click Add Devotee"
$ws.Range("D2").Value = "The application to be added successfully"
$ws.Range("E2").Value = "The application to be added successfully- test failed"
Set-TextValue $ws.Range("F2") "1"
$ws.Range("G2").Value = "1 - Critical"

# Rows 3-9: only column C (Test Procedure) changes
for ($r = 3; $r -le 9; $r++) {
    $ws.Range("C$r").Value = "1) Click Add Devotee"
}

# Row 10
Set-TextValue $ws.Range("A10") "13083"
$ws.Range("B10").Value = "Scenario2"
$ws.Range("C10").Value = "1) Click Add Devotee"
$ws.Range("D10").Value = "The application to be added successfully"
$ws.Range("E10").Value = "The application to be added successfully- test failed"
Set-TextValue $ws.Range("F10") "2"
$ws.Range("G10").Value = "2 - High"

# Rows 11-18: only column C (Test Procedure) changes
for ($r = 11; $r -le 18; $r++) {
    $ws.Range("C$r").Value = "1) Click Add Devotee"
}

$excel.CutCopyMode = $false
